# Diario de fusao - append May 14-19 daily model totals, rows 812:877
# Rows 812-822 (14/05/2025) were pasted in originally as text; retype them as
# numbers. New rows 823-866 (14-17/05/2025) are entered as numbers; the final
# batch for 19/05/2025 (rows 867-877) keeps the text formatting used so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entries = @(
    @{ Row = 812; Date = "14/05/2025"; Model = "2811"; Value = "650"; AsText = $false },
    @{ Row = 813; Date = "14/05/2025"; Model = "3581"; Value = "457"; AsText = $false },
    @{ Row = 814; Date = "14/05/2025"; Model = "2651"; Value = "135"; AsText = $false },
    @{ Row = 815; Date = "14/05/2025"; Model = "3410"; Value = "333"; AsText = $false },
    @{ Row = 816; Date = "14/05/2025"; Model = "3760"; Value = "273"; AsText = $false },
    @{ Row = 817; Date = "14/05/2025"; Model = "3631"; Value = "201"; AsText = $false },
    @{ Row = 818; Date = "14/05/2025"; Model = "3630"; Value = "5"; AsText = $false },
    @{ Row = 819; Date = "14/05/2025"; Model = "3761"; Value = "10"; AsText = $false },
    @{ Row = 820; Date = "14/05/2025"; Model = "3591"; Value = "51"; AsText = $false },
    @{ Row = 821; Date = "14/05/2025"; Model = "3752"; Value = "69"; AsText = $false },
    @{ Row = 822; Date = "14/05/2025"; Model = "3771"; Value = "58"; AsText = $false },
    @{ Row = 823; Date = "14/05/2025"; Model = "2811"; Value = "650"; AsText = $false },
    @{ Row = 824; Date = "14/05/2025"; Model = "3581"; Value = "457"; AsText = $false },
    @{ Row = 825; Date = "14/05/2025"; Model = "2651"; Value = "135"; AsText = $false },
    @{ Row = 826; Date = "14/05/2025"; Model = "3410"; Value = "333"; AsText = $false },
    @{ Row = 827; Date = "14/05/2025"; Model = "3760"; Value = "273"; AsText = $false },
    @{ Row = 828; Date = "14/05/2025"; Model = "3631"; Value = "201"; AsText = $false },
    @{ Row = 829; Date = "14/05/2025"; Model = "3630"; Value = "5"; AsText = $false },
    @{ Row = 830; Date = "14/05/2025"; Model = "3761"; Value = "10"; AsText = $false },
    @{ Row = 831; Date = "14/05/2025"; Model = "3591"; Value = "51"; AsText = $false },
    @{ Row = 832; Date = "14/05/2025"; Model = "3752"; Value = "69"; AsText = $false },
    @{ Row = 833; Date = "14/05/2025"; Model = "3771"; Value = "58"; AsText = $false },
    @{ Row = 834; Date = "15/05/2025"; Model = "3410"; Value = "311"; AsText = $false },
    @{ Row = 835; Date = "15/05/2025"; Model = "3631"; Value = "222"; AsText = $false },
    @{ Row = 836; Date = "15/05/2025"; Model = "3760"; Value = "229"; AsText = $false },
    @{ Row = 837; Date = "15/05/2025"; Model = "3581"; Value = "467"; AsText = $false },
    @{ Row = 838; Date = "15/05/2025"; Model = "2651"; Value = "116"; AsText = $false },
    @{ Row = 839; Date = "15/05/2025"; Model = "2811"; Value = "601"; AsText = $false },
    @{ Row = 840; Date = "15/05/2025"; Model = "3630"; Value = "5"; AsText = $false },
    @{ Row = 841; Date = "15/05/2025"; Model = "3761"; Value = "10"; AsText = $false },
    @{ Row = 842; Date = "15/05/2025"; Model = "3591"; Value = "119"; AsText = $false },
    @{ Row = 843; Date = "15/05/2025"; Model = "3771"; Value = "52"; AsText = $false },
    @{ Row = 844; Date = "15/05/2025"; Model = "3752"; Value = "34"; AsText = $false },
    @{ Row = 845; Date = "16/05/2025"; Model = "3410"; Value = "332"; AsText = $false },
    @{ Row = 846; Date = "16/05/2025"; Model = "2651"; Value = "143"; AsText = $false },
    @{ Row = 847; Date = "16/05/2025"; Model = "3761"; Value = "46"; AsText = $false },
    @{ Row = 848; Date = "16/05/2025"; Model = "3760"; Value = "193"; AsText = $false },
    @{ Row = 849; Date = "16/05/2025"; Model = "3581"; Value = "385"; AsText = $false },
    @{ Row = 850; Date = "16/05/2025"; Model = "3630"; Value = "5"; AsText = $false },
    @{ Row = 851; Date = "16/05/2025"; Model = "2811"; Value = "660"; AsText = $false },
    @{ Row = 852; Date = "16/05/2025"; Model = "3631"; Value = "196"; AsText = $false },
    @{ Row = 853; Date = "16/05/2025"; Model = "3591"; Value = "91"; AsText = $false },
    @{ Row = 854; Date = "16/05/2025"; Model = "3771"; Value = "88"; AsText = $false },
    @{ Row = 855; Date = "16/05/2025"; Model = "3752"; Value = "87"; AsText = $false },
    @{ Row = 856; Date = "17/05/2025"; Model = "3760"; Value = "117"; AsText = $false },
    @{ Row = 857; Date = "17/05/2025"; Model = "3410"; Value = "289"; AsText = $false },
    @{ Row = 858; Date = "17/05/2025"; Model = "2811"; Value = "618"; AsText = $false },
    @{ Row = 859; Date = "17/05/2025"; Model = "3581"; Value = "337"; AsText = $false },
    @{ Row = 860; Date = "17/05/2025"; Model = "3761"; Value = "88"; AsText = $false },
    @{ Row = 861; Date = "17/05/2025"; Model = "2651"; Value = "133"; AsText = $false },
    @{ Row = 862; Date = "17/05/2025"; Model = "3630"; Value = "2"; AsText = $false },
    @{ Row = 863; Date = "17/05/2025"; Model = "3631"; Value = "167"; AsText = $false },
    @{ Row = 864; Date = "17/05/2025"; Model = "3591"; Value = "45"; AsText = $false },
    @{ Row = 865; Date = "17/05/2025"; Model = "3752"; Value = "23"; AsText = $false },
    @{ Row = 866; Date = "17/05/2025"; Model = "3771"; Value = "62"; AsText = $false },
    @{ Row = 867; Date = "19/05/2025"; Model = "2651"; Value = "93"; AsText = $true },
    @{ Row = 868; Date = "19/05/2025"; Model = "2811"; Value = "724"; AsText = $true },
    @{ Row = 869; Date = "19/05/2025"; Model = "3410"; Value = "367"; AsText = $true },
    @{ Row = 870; Date = "19/05/2025"; Model = "3581"; Value = "355"; AsText = $true },
    @{ Row = 871; Date = "19/05/2025"; Model = "3761"; Value = "240"; AsText = $true },
    @{ Row = 872; Date = "19/05/2025"; Model = "3631"; Value = "214"; AsText = $true },
    @{ Row = 873; Date = "19/05/2025"; Model = "3630"; Value = "3"; AsText = $true },
    @{ Row = 874; Date = "19/05/2025"; Model = "3760"; Value = "36"; AsText = $true },
    @{ Row = 875; Date = "19/05/2025"; Model = "3591"; Value = "88"; AsText = $true },
    @{ Row = 876; Date = "19/05/2025"; Model = "3771"; Value = "16"; AsText = $true },
    @{ Row = 877; Date = "19/05/2025"; Model = "3752"; Value = "14"; AsText = $true }
)

foreach ($e in $entries) {
    $ws.Range("A$($e.Row)").Value = $e.Date
    if ($e.AsText) {
        $ws.Range("B$($e.Row)").NumberFormat = "@"
        $ws.Range("B$($e.Row)").Value = $e.Model
        $ws.Range("B$($e.Row)").Style = "Normal"
        $ws.Range("C$($e.Row)").NumberFormat = "@"
        $ws.Range("C$($e.Row)").Value = $e.Value
        $ws.Range("C$($e.Row)").Style = "Normal"
    } else {
        $ws.Range("B$($e.Row)").Value = [double]$e.Model
        $ws.Range("C$($e.Row)").Value = [double]$e.Value
    }
}
